$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text formatting is preserved for Price (D) and Volume (E) columns
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.390.79"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "1.881.11"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "0.7164"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").Value = "243.45"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "0.07946"
$ws.Range("E8").Value = "  -1.23%  "
$ws.Range("D9").Value = "0.3138"
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("E10").Value = "  -1.63%  "
$ws.Range("D11").Value = "0.08083"
$ws.Range("E11").Value = "  -3.33%  "
$ws.Range("D12").Value = "1.917.16"
$ws.Range("E12").Value = "  +2.08%  "
$ws.Range("E13").Value = "  +3.49%  "
$ws.Range("D14").Value = "5.203"
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("D15").Value = "0.7072"
$ws.Range("E15").Value = "  -1.77%  "
$ws.Range("D16").Value = "6.382"
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("D17").Value = "0.000008403"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "29.524.71"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").Value = "252.42"
$ws.Range("E19").Value = "  +4.75%  "
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("D21").Value = "2.137.00"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "7.680"
$ws.Range("E23").Value = "  -1.66%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "0.1575"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").Value = "161.87"
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("D28").Value = "18.93"
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "4.409"
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("D31").Value = "4.311"
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("D32").Value = "1.236"
$ws.Range("E32").Value = "  +3.39%  "
$ws.Range("D33").Value = "0.05298"
$ws.Range("E33").Value = "  -1.50%  "
$ws.Range("D34").Value = "1.932"
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("D35").Value = "0.7555"
$ws.Range("E35").Value = "  +0.60%  "
$ws.Range("E36").Value = "  -0.42%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").Value = "1.285.12"
$ws.Range("E38").Value = "  -0.71%  "
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("D40").Value = "2.765"
$ws.Range("E40").Value = "  +0.96%  "
$ws.Range("D41").Value = "6.392"
$ws.Range("E41").Value = "  -2.94%  "
$ws.Range("D42").Value = "0.9059"
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("D43").Value = "111.63"
$ws.Range("E43").Value = "  +0.84%  "
$ws.Range("D44").Value = "73.91"
$ws.Range("E44").Value = "  +1.13%  "
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("D47").Value = "2.039.62"
$ws.Range("E47").Value = "  +0.77%  "
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").Value = "9.501"
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("E51").Value = "  -0.56%  "
